$wb = $excel.ActiveWorkbook

# --- "Change History" sheet: append a new history row (CRE21-005) -------
$wsContent = $wb.Worksheets.Item("Content")
$wsCH = $wb.Worksheets.Item("Change History")

# Seed the new row's formatting by copying from existing rows that already
# carry the font/number-format combinations we need, then retarget the
# vertical alignment to "top" (matching the template's newer history rows)
# and finally overwrite with the real values.
$wsCH.Range("B6").Copy($wsCH.Range("A7"))
$wsCH.Range("B6").Copy($wsCH.Range("B7"))
$wsCH.Range("C4").Copy($wsCH.Range("C7"))
$wsCH.Range("D6").Copy($wsCH.Range("D7"))

$wsCH.Range("A7:D7").VerticalAlignment = -4160

$wsCH.Range("A7").Value = 4
$wsCH.Range("B7").Value = "CRE21-005"
$wsCH.Range("C7").Value = "To exclude claim transactions related to COVID-19 vaccine subsidies in PPC0002, PPC0003, eHS(S)U0002"
$wsCH.Range("D7").Value = 44351

$wsCH.Rows.Item(7).RowHeight = 28.5

# Move the "Change History" sheet's remembered selection to A3 (as in the
# updated template) without disturbing which tab is actually active.
[void]$wsCH.Activate()
[void]$wsCH.Range("A3").Select()
[void]$wsContent.Activate()
